$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the risk value for the first FIXED RATIO block (D6): 0.03% -> 0.15%
$ws.Range("D6").Value = 0.0015

# Add the new "bet amount" columns (H/I) for both FIXED RATIO blocks.
# Row 6 block
$ws.Range("H6").Formula = "=D6*H2"
$ws.Range("H6").NumberFormat = "0.00%"
$ws.Range("I6").Formula = "=E6*H6"
$ws.Range("I6").NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"

# Row 10 block
$ws.Range("G10").Formula = "=E10*D10"
$ws.Range("G10").NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"
$ws.Range("H10").Formula = "=D10*H2"
$ws.Range("H10").NumberFormat = "0.00%"
$ws.Range("I10").Formula = "=E10*H10"
$ws.Range("I10").NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"

# Leave the selection where the author ended up after making these edits.
[void]$ws.Range("O6").Select()
